# ---------------------------------------------------------------------------
# Adds phone/email/address to Users, a category column (+ a maintenance row)
# to Bookings, and three brand-new sheets: OperatingHours, ClosedDates,
# Settings.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Users sheet — add phone / email / address columns (C, D, E)
# ---------------------------------------------------------------------------
$users = $wb.Worksheets.Item("Users")

# Header row: clone the existing bold/centered/bordered header style from B1.
$users.Range("B1").Copy()
$users.Range("C1:E1").PasteSpecial(-4122)
$users.Range("C1").Value = "phone"
$users.Range("D1").Value = "email"
$users.Range("E1").Value = "address"

$firstNames = @("alex","ben","charlie","daniel","ethan","finn","george","harry","isaac","jack","kai","liam","mason","noah","oscar","peter","quentin","riley","samuel","thomas","umar","victor","william","yusuf","zachary")

for ($i = 0; $i -lt $firstNames.Count; $i++) {
    $row = $i + 2
    $phoneSuffix = "{0:D3}" -f $i
    $streetNum = $i + 11
    # Leading '0' prefix forces these to stay text (otherwise Excel would
    # coerce the numeric-looking phone number and drop the leading zeros).
    $users.Range("C$row").Value = "'0400100$phoneSuffix"
    $users.Range("D$row").Value = "$($firstNames[$i])@example.com"
    $users.Range("E$row").Value = "$streetNum Workshop Rd, Newcastle NSW 2300"
}

# ---------------------------------------------------------------------------
# 2. Bookings sheet — add category column (G) + a new Maintenance row (7)
# ---------------------------------------------------------------------------
$bookings = $wb.Worksheets.Item("Bookings")

$bookings.Range("F1").Copy()
$bookings.Range("G1").PasteSpecial(-4122)
$bookings.Range("G1").Value = "category"

for ($row = 2; $row -le 6; $row++) {
    $bookings.Range("G$row").Value = "Usage"
}

# Row 7: new maintenance booking. Clone the date-time style from row 6's
# start/end columns first, then write the values.
$bookings.Range("D6:E6").Copy()
$bookings.Range("D7").PasteSpecial(-4122)

$bookings.Range("A7").Value = 6
$bookings.Range("B7").Value = 0
$bookings.Range("C7").Value = 3
$bookings.Range("D7").Value = 45914.54166666666
$bookings.Range("E7").Value = 45914.60416666666
$bookings.Range("F7").Value = "Confirmed"
$bookings.Range("G7").Value = "Maintenance"

# ---------------------------------------------------------------------------
# 3. New sheets — appended after ServiceLog (the current last sheet)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$hours = $wb.Worksheets.Add($null, $lastSheet)
$hours.Name = "OperatingHours"

$closed = $wb.Worksheets.Add($null, $hours)
$closed.Name = "ClosedDates"

$settings = $wb.Worksheets.Add($null, $closed)
$settings.Name = "Settings"

# --- OperatingHours ----------------------------------------------------------
$users.Range("B1").Copy()
$hours.Range("A1:E1").PasteSpecial(-4122)
$hours.Range("A1").Value = "weekday"
$hours.Range("B1").Value = "name"
$hours.Range("C1").Value = "is_open"
$hours.Range("D1").Value = "open"
$hours.Range("E1").Value = "close"

$dayNames = @("Mon","Tue","Wed","Thu","Fri","Sat","Sun")
$dayOpen = @($true,$true,$true,$true,$true,$true,$false)
$dayStart = @("08:00","08:00","08:00","08:00","08:00","09:00","00:00")
$dayEnd = @("17:00","17:00","17:00","17:00","17:00","13:00","00:00")

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $hours.Range("A$row").Value = $i
    $hours.Range("B$row").Value = $dayNames[$i]
    $hours.Range("C$row").Value = $dayOpen[$i]
    $hours.Range("D$row").Value = $dayStart[$i]
    $hours.Range("E$row").Value = $dayEnd[$i]
}

# --- ClosedDates --------------------------------------------------------------
$users.Range("B1").Copy()
$closed.Range("A1:B1").PasteSpecial(-4122)
$closed.Range("A1").Value = "date"
$closed.Range("B1").Value = "reason"

# --- Settings -----------------------------------------------------------------
$users.Range("B1").Copy()
$settings.Range("A1:B1").PasteSpecial(-4122)
$settings.Range("A1").Value = "key"
$settings.Range("B1").Value = "value"

$settings.Range("A2").Value = "admin_password"
$settings.Range("B2").Value = "nesnob2025"
$settings.Range("A3").Value = "show_contact_on_bookings"
# Leading "'" forces literal text "true" instead of an auto-coerced boolean.
$settings.Range("B3").Value = "'true"

Write-Output "edit complete"
